$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 118; existing rows 118-138 shift down to 119-139.
$ws.Rows.Item(118).Insert()

$newRow = 118

$ws.Cells.Item($newRow, 1).Value = 4
$ws.Cells.Item($newRow, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($newRow, 3).Value = "Los Lagos"
$ws.Cells.Item($newRow, 4).Value = 44476
$ws.Cells.Item($newRow, 5).Value = 10
$ws.Cells.Item($newRow, 6).Value = 100112017
$ws.Cells.Item($newRow, 7).Value = "Apio"
$ws.Cells.Item($newRow, 8).Value = "Americana (o)"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 35
$ws.Cells.Item($newRow, 11).Value = 11000
$ws.Cells.Item($newRow, 12).Value = 12000
$ws.Cells.Item($newRow, 13).Value = 11429
$ws.Cells.Item($newRow, 14).Value = '$/docena de matas'
$ws.Cells.Item($newRow, 15).Value = "Región de Coquimbo"
$ws.Cells.Item($newRow, 16).Value = 1905
$ws.Cells.Item($newRow, 17).Value = 6
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
